# Refresh cryptocurrency Price (D) and Volume(1h) (E) columns with the
# latest values scraped by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.575.22'
$ws.Range("E2").Value = '  +1.74%  '
$ws.Range("D3").Value = '3.020.76'
$ws.Range("E3").Value = '  +1.70%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = "'596.66"
$ws.Range("E5").Value = '  +1.52%  '
$ws.Range("D6").Value = "'150.62"
$ws.Range("E6").Value = '  +6.64%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '3.018.36'
$ws.Range("E8").Value = '  +1.70%  '
$ws.Range("D9").Value = "'0.519"
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").Value = "'6.40"
$ws.Range("E10").Value = '  +11.44%  '
$ws.Range("E11").Value = '  +4.85%  '
$ws.Range("D12").Value = "'0.460"
$ws.Range("E12").Value = '  +1.07%  '
$ws.Range("D14").Value = "'34.61"
$ws.Range("E14").Value = '  +2.16%  '
$ws.Range("D16").Value = '3.520.38'
$ws.Range("E16").Value = '  +1.56%  '
$ws.Range("D17").Value = '62.514.36'
$ws.Range("E17").Value = '  +1.64%  '
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("D19").Value = '3.021.59'
$ws.Range("E19").Value = '  +1.75%  '
$ws.Range("D20").Value = "'449.74"
$ws.Range("E20").Value = '  +0.23%  '
$ws.Range("E21").Value = '  +2.72%  '
$ws.Range("D22").Value = "'0.691"
$ws.Range("E22").Value = '  +1.58%  '
$ws.Range("D23").Value = "'7.47"
$ws.Range("E23").Value = '  +2.10%  '
$ws.Range("D24").Value = "'82.39"
$ws.Range("E24").Value = '  +1.59%  '
$ws.Range("D25").Value = "'10.90"
$ws.Range("E25").Value = '  +11.56%  '
$ws.Range("E26").Value = '  +5.13%  '
$ws.Range("D27").Value = "'12.10"
$ws.Range("E27").Value = '  +0.27%  '
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("E29").Value = '  +3.21%  '
$ws.Range("D30").Value = "'7.36"
$ws.Range("E30").Value = '  +7.97%  '
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("E32").Value = '  +4.87%  '
$ws.Range("D33").Value = "'27.54"
$ws.Range("E33").Value = '  +2.17%  '
$ws.Range("E34").Value = '  +3.16%  '
$ws.Range("D35").Value = '0.0₃0856'
$ws.Range("E35").Value = '  +10.99%  '
$ws.Range("E36").Value = '  +1.53%  '
$ws.Range("E37").Value = '  +2.90%  '
$ws.Range("D38").Value = "'3.08"
$ws.Range("E38").Value = '  +11.66%  '
$ws.Range("D39").Value = "'2.08"
$ws.Range("E39").Value = '  +0.94%  '
$ws.Range("D40").Value = "'50.17"
$ws.Range("E40").Value = '  +0.18%  '
$ws.Range("D41").Value = "'9.04"
$ws.Range("E41").Value = '  -0.66%  '
$ws.Range("E42").Value = '  +4.63%  '
$ws.Range("D43").Value = "'0.288"
$ws.Range("E43").Value = '  +9.97%  '
$ws.Range("D44").Value = "'40.55"
$ws.Range("E44").Value = '  +10.26%  '
$ws.Range("D45").Value = "'393.54"
$ws.Range("E45").Value = '  +2.34%  '
$ws.Range("E46").Value = '  +0.84%  '
$ws.Range("D47").Value = '2.738.04'
$ws.Range("E47").Value = '  +1.08%  '
$ws.Range("D48").Value = "'132.61"
$ws.Range("E48").Value = '  +2.14%  '
$ws.Range("D50").Value = "'2.18"
$ws.Range("E50").Value = '  +1.48%  '
$ws.Range("E51").Value = '  +0.14%  '
